# Apply "6 month coverage from 2026 in scenario 3a" edit.
#
# Both sheets currently have one column per whole year (2018..2040) in the
# header row. The edit doubles the time resolution to half-year steps
# (2018, 2018.5, 2019, 2019.5, ..., 2040) by inserting a new blank column
# immediately after each existing year column, then (re)writing every
# header/data cell explicitly from its target (year -> column) mapping, so
# the exact pre-insert cell shuffling doesn't matter:
#   - header row gets every half-year label
#   - scenario rows that start in 2026 now report a value every 6 months
#     (instead of every 2 years)
#   - the pre-2026 row reports a value every year (instead of every 2 years)
#   - the MarketShare rows likewise go from annual to semi-annual

$wb = $excel.ActiveWorkbook

function Insert-HalfYearColumns($ws, $firstCol, $lastCol) {
    # Walk from the last original year-column back to the first, inserting
    # one blank column right after each, so the sheet ends up with twice
    # as many (now half-year-spaced) columns and the correct dimension.
    for ($col = $lastCol; $col -ge $firstCol; $col--) {
        $ws.Cells.Item(1, $col + 1).EntireColumn.Insert()
    }
}

# ---------------------------------------------------------------------
# Sheet 1: "Platform Coverage"
# Year columns before edit: H(8) .. AD(30) => 2018 .. 2040
# After edit (half-year step): H(8) .. AZ(52) => 2018 .. 2040
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Platform Coverage")

Insert-HalfYearColumns $ws1 8 30

# Header row: every column from H to AZ gets year 2018 + 0.5*(col-8).
for ($col = 8; $col -le 52; $col++) {
    $ws1.Cells.Item(1, $col).Value = 2018.0 + 0.5 * ($col - 8)
}

# Row 2 (0.6 coverage): now every whole year 2018-2025 (cols 8,10,...,22).
for ($col = 8; $col -le 22; $col = $col + 2) {
    $ws1.Cells.Item(2, $col).Value = 0.6
}

# Rows 3-5 (coverage from 2026 onward): now every half year 2026-2040
# (cols 24..52, i.e. every single column in that span).
for ($col = 24; $col -le 52; $col++) {
    $ws1.Cells.Item(3, $col).Value = 0.8
    $ws1.Cells.Item(4, $col).Value = 0.5
    $ws1.Cells.Item(5, $col).Value = 0.5
}

# View bookkeeping to match the saved workbook.
$ws1.Range("AZ17").Select()
$ws1.Application.ActiveWindow.ScrollColumn = $ws1.Range("AG1").Column

# ---------------------------------------------------------------------
# Sheet 2: "MarketShare"
# Year columns before edit: D(4) .. Z(26) => 2018 .. 2040
# After edit (half-year step): D(4) .. AV(48) => 2018 .. 2040
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MarketShare")

Insert-HalfYearColumns $ws2 4 26

# Header row: every column from D to AV gets year 2018 + 0.5*(col-4).
for ($col = 4; $col -le 48; $col++) {
    $ws2.Cells.Item(1, $col).Value = 2018.0 + 0.5 * ($col - 4)
}

# Row 2 (value 1): now every half year 2026-2040 (cols 20..48).
for ($col = 20; $col -le 48; $col++) {
    $ws2.Cells.Item(2, $col).Value = 1
}

# Row 3 (value 1): now every half year 2018-2025.5 (cols 4..19).
for ($col = 4; $col -le 19; $col++) {
    $ws2.Cells.Item(3, $col).Value = 1
}

# View bookkeeping to match the saved workbook.
$ws2.Range("Q3").Select()
$ws2.Application.ActiveWindow.ScrollColumn = $ws2.Range("H1").Column
